# EMPANADA Data Files Key.xlsx - "Took more data (days 8-11)"
#
# Adds 40 new rows (80-119) for Day8..Day11 Earth-gravity entries, with
# Original File Name values MVI_0001..MVI_0029, MVI_0041, MVI_0030..MVI_0038,
# MVI_0040 (in that exact order), a shared "=Day..." concatenation formula in
# column B, and restyles the previously "italic/serif-ish" (style 5) rows
# 61-79 (and a couple of stray numeric cells) back to the normal style 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New rows 80-119: values (Original File Name, Day, Gravity, Speed)
# ---------------------------------------------------------------------

$names = @(
  "MVI_0001.MOV","MVI_0002.MOV","MVI_0003.MOV","MVI_0004.MOV","MVI_0005.MOV",
  "MVI_0006.MOV","MVI_0007.MOV","MVI_0008.MOV","MVI_0009.MOV","MVI_0010.MOV",
  "MVI_0011.MOV","MVI_0012.MOV","MVI_0013.MOV","MVI_0014.MOV","MVI_0015.MOV",
  "MVI_0016.MOV","MVI_0017.MOV","MVI_0018.MOV","MVI_0019.MOV","MVI_0020.MOV",
  "MVI_0021.MOV","MVI_0022.MOV","MVI_0023.MOV","MVI_0024.MOV","MVI_0025.MOV",
  "MVI_0026.MOV","MVI_0027.MOV","MVI_0028.MOV","MVI_0029.MOV","MVI_0041.MOV",
  "MVI_0030.MOV","MVI_0031.MOV","MVI_0032.MOV","MVI_0033.MOV","MVI_0034.MOV",
  "MVI_0035.MOV","MVI_0036.MOV","MVI_0037.MOV","MVI_0038.MOV","MVI_0040.MOV"
)

$days   = @(8,8,8,8,8,8,8,8,8,8, 9,9,9,9,9,9,9,9,9,9, 10,10,10,10,10,10,10,10,10,10, 11,11,11,11,11,11,11,11,11,11)
$speeds = @(1,2,3,4,5,6,7,8,9,10, 1,2,3,4,5,6,7,8,9,10, 1,2,3,4,5,6,7,8,9,10, 1,2,3,4,5,6,7,8,9,10)

for ($i = 0; $i -lt 40; $i++) {
    $r = 80 + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $days[$i]
    $ws.Cells.Item($r, 4).Value = "Earth"
    $ws.Cells.Item($r, 5).Value = $speeds[$i]
}

# Column B: one shared "=Day"&C&"-"&D&"-"&E&"mms.mov" formula across B80:B119
$ws.Range("B80:B119").Formula = '="Day"&C80&"-"&D80&"-"&E80&"mms.mov"'

# ---------------------------------------------------------------------
# 2) Copy cell formatting (fonts/styles only) from existing cells that
#    already carry the styles this new block needs, onto the new cells.
#    Done BEFORE re-styling rows 61-79 below, since those rows are the
#    source of the "style 5" look.
# ---------------------------------------------------------------------

# Column A: style 5 (default/Calibri font) for most rows, using A61 as the
# donor (A61 is itself being restyled to style 2 later in this script).
$ws.Range("A61").Copy() | Out-Null
$ws.Range("A81").PasteSpecial(-4122) | Out-Null
$ws.Range("A83:A108").PasteSpecial(-4122) | Out-Null
$ws.Range("A110:A118").PasteSpecial(-4122) | Out-Null

# Column A: style 2 (Arial) for the exceptions - A55 already carries it.
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A80").PasteSpecial(-4122) | Out-Null
$ws.Range("A82").PasteSpecial(-4122) | Out-Null
$ws.Range("A109").PasteSpecial(-4122) | Out-Null
$ws.Range("A119").PasteSpecial(-4122) | Out-Null

# Column B: style 3 (Arial, formula cell) for the whole block - B33 already
# carries it (it is itself a shared "Day..." formula cell).
$ws.Range("B33").Copy() | Out-Null
$ws.Range("B80:B119").PasteSpecial(-4122) | Out-Null

# Column C: style 5 for C80/C110, donor F76 (numeric, style 5, untouched by
# this edit) ; style 2 for the rest, donor E56 (numeric, style 2).
$ws.Range("F76").Copy() | Out-Null
$ws.Range("C80").PasteSpecial(-4122) | Out-Null
$ws.Range("C110").PasteSpecial(-4122) | Out-Null

$ws.Range("E56").Copy() | Out-Null
$ws.Range("C81:C109").PasteSpecial(-4122) | Out-Null
$ws.Range("C111:C119").PasteSpecial(-4122) | Out-Null

# Columns D & E: style 2 throughout, donor D2 / E56.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D80:D119").PasteSpecial(-4122) | Out-Null

$ws.Range("E56").Copy() | Out-Null
$ws.Range("E80:E119").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Restyle existing cells: rows 61-79 (plus F56 and G72) move from the
#    "style 5" look to the normal "style 2" look used everywhere else.
# ---------------------------------------------------------------------

$ws.Range("E56").Copy() | Out-Null
$ws.Range("F56").PasteSpecial(-4122) | Out-Null
$ws.Range("G72").PasteSpecial(-4122) | Out-Null

$ws.Range("A55").Copy() | Out-Null
$ws.Range("A61:B79").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
